# update detail satuan sampah dan add fitur cancel setoran
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old extra waste entries (rows 4-13), shifting cells up so only
# the header row plus two data rows (A1:D3) remain.
$ws.Range("A4:D13").Delete(-4162)  # xlShiftUp

# Update the "Tanggal Update" (date) cells for the remaining rows.
# Force these to stay plain text (not get auto-converted to a date serial
# number) by temporarily marking the cell as text-formatted, then clearing
# the formatting back off again once the text value has been stored.
$dateCells = $ws.Range("B2:B3")
$dateCells.NumberFormat = "@"
$dateCells.Value = "2024-03-18"
$dateCells.ClearFormats()

# Update the "Jenis Limbah" values.
$ws.Range("C2").Value = "cair"
$ws.Range("C3").Value = "padat"

# Update the price for row 3.
$ws.Range("D3").Value = 15000
